# Updated symbol list on Mon Dec 19 09:30:02 UTC 2022 with GitHub Actions
#
# Refresh the crypto price snapshot: most rows just get a small price
# tick, while a handful of rows (17-25 and 41,43) are re-ranked so the
# Coin/Link/Volume columns shift to the neighbouring coin.
#
# NOTE: column D values are numeric-looking strings that must stay text
# (the sheet stores them as inline/shared strings, not numbers), so they
# are written with a leading apostrophe to force Excel to keep them as
# text instead of auto-converting to a float.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumber($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

# --- simple price refreshes (rows 2-16) ---
Set-TextNumber "D2"  "247.60"
Set-TextNumber "D3"  "21.77"
Set-TextNumber "D4"  "5.462"
Set-TextNumber "D5"  "0.05693"
Set-TextNumber "D6"  "3.381"
Set-TextNumber "D7"  "0.8064"
Set-TextNumber "D8"  "1.036"
Set-TextNumber "D9"  "0.1498"
Set-TextNumber "D10" "0.07307"
Set-TextNumber "D11" "0.03169"
Set-TextNumber "D12" "0.02939"
Set-TextNumber "D13" "0.09285"
Set-TextNumber "D14" "0.001646"
Set-TextNumber "D15" "3.212"
Set-TextNumber "D16" "0.04714"

# --- rows 17-25: coin ranking rotates by one position, with "One" wrapping
#     from row 17 down to row 25 ---
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextNumber "D17" "0.006343"
$ws.Range("E17").Value = "16TigerCashTCH"

$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextNumber "D18" "0.005047"
$ws.Range("E18").Value = "17HotbitTokenHTB"

$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextNumber "D19" "0.001046"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "NitroEx"
$ws.Range("C20").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextNumber "D20" "0.0001501"
$ws.Range("E20").Value = "19NitroExNTX"

$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextNumber "D21" "0.0003201"
$ws.Range("E21").Value = "20UpBotsUBXT"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextNumber "D22" "3.776"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "KuCoinToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextNumber "D23" "6.426"
$ws.Range("E23").Value = "22KuCoinTokenKCS"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextNumber "D24" "2.094"
$ws.Range("E24").Value = "23BTSETokenBTSE"

$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextNumber "D25" "0.01160"
$ws.Range("E25").Value = "24OneONEBestin24h"

# --- row 40: simple price refresh ---
Set-TextNumber "D40" "0.04111"

# --- rows 41 & 43 swap places (row 42 stays put) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextNumber "D41" "0.006943"
$ws.Range("E41").Value = "40KickTokenKICK"

Set-TextNumber "D42" "0.003501"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextNumber "D43" "0.1043"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- remaining simple price refreshes ---
Set-TextNumber "D44" "0.008112"
Set-TextNumber "D45" "0.00005832"
Set-TextNumber "D47" "0.0005800"
Set-TextNumber "D49" "0.009694"
